# "Few fixes for new release"
#  - fix the "Rolland_Garros" typo -> "Roland_Garros" on every sheet's header (C1)
#  - leave the last selection on each sheet the way the author left it, and
#    make "Novak Djokovic" (3rd sheet) the active/selected tab on reopen

$wb = $excel.ActiveWorkbook

$wsFederer  = $wb.Worksheets.Item(1)   # Roger Federer
$wsNadal    = $wb.Worksheets.Item(2)   # Rafael Nadal
$wsDjokovic = $wb.Worksheets.Item(3)   # Novak Djokovic

# Fix the misspelled "Rolland_Garros" column header on every sheet.
$wsFederer.Range("C1").Value  = "Roland_Garros"
$wsNadal.Range("C1").Value    = "Roland_Garros"
$wsDjokovic.Range("C1").Value = "Roland_Garros"

# Restore each sheet's own last-used selection.
$wsFederer.Activate()
$wsFederer.Range("H6").Select()

$wsNadal.Activate()
$wsNadal.Range("E7").Select()

$wsDjokovic.Activate()
$wsDjokovic.Range("H7").Select()
